# Applies the "Add files via upload" commit:
#  1. Update the subtitle on the title slide (slide 1) to append the
#     student id suffix.
#  2. Append a new blank-layout slide at the end of the deck containing a
#     single "QUESTION ^_^" textbox.

$p = $ppt.ActivePresentation

# --- 1. Title slide subtitle text -----------------------------------------
$titleSlide = $p.Slides.Item(1)
$subtitleShape = $titleSlide.Shapes.Item(2)
$subtitleRun = $subtitleShape.TextFrame.TextRange.Runs(1)
$subtitleRun.Text = "Presented by Trong-Binh Nguyen-202388548"

# --- 2. New "QUESTION ^_^" slide -------------------------------------------
# ppLayoutBlank = 12
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)

# EMU -> points (1 pt = 12700 EMU)
$left   = 4168877 / 12700
$top    = 2721114 / 12700
$width  = 6676103 / 12700
$height = 707886 / 12700

$textBox = $newSlide.Shapes.AddTextbox(1, $left, $top, $width, $height)
$textBox.Fill.Visible = $false
$textBox.TextFrame.WordWrap = $true
$textBox.TextFrame.AutoSize = 1

$questionRange = $textBox.TextFrame.TextRange
$questionRange.Text = "QUESTION ^_^"
$questionRange.Font.Size = 40
$questionRange.Font.Name = "Times New Roman"
$questionRange.Font.NameComplexScript = "Times New Roman"
